$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as a number by Excel's type
# coercion (e.g. "1.000", "0.7164"). Force them to stay literal text by
# switching the cell to Text format before the write, then restoring the
# cell style afterwards so no stray formatting is left behind.
$ws.Range("D2").Value = '29.451.79'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = '1.875.92'
$ws.Range("E3").Value = '  +0.97%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.7164'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.97%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '241.41'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("E7").Value = '  +0.01%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.07909'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.44%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.3095'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +2.05%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '25.45'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +4.33%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.08262'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.10%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.7273'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.31%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.869.49'
$ws.Range("E13").Value = '  +0.23%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.268'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.17%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '91.25'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.89%  '
$ws.Range("D16").Value = '29.425.46'
$ws.Range("E16").Value = '  +0.38%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '5.895'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.73%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '245.53'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +3.38%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000007843'
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.28'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = '2.128.00'
$ws.Range("E21").Value = '  +0.25%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '8.053'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +6.62%  '
$ws.Range("E23").Value = '  +0.07%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E25").Value = '  +14.46%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '162.64'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.20%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.048'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.73%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '18.30'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.12%  '
$ws.Range("E29").Value = '  -3.07%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.492'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.06%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.410'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.19%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.103'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.33%  '
$ws.Range("E33").Value = '  +0.39%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.947'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.80%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.197'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.84%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.7250'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.70%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.675'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.06%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01868'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.18%  '
$ws.Range("D39").Value = '1.199.57'
$ws.Range("E39").Value = '  +4.77%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.690'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.21%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.9084'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.49%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '6.149'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +3.27%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '73.11'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.71%  '
$ws.Range("E44").Value = '  +0.02%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '102.12'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.41%  '
$ws.Range("D46").Value = '2.020.54'
$ws.Range("E46").Value = '  +0.03%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.5287'
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.794'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.39%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.917'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +8.09%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '9.298'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.39%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.4316'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.50%  '
